$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 238, shifting the existing
# rows 238-255 down to 239-256 (preserves formatting of the row below,
# e.g. the date-formatted style on column D).
$ws.Rows.Item(238).Insert()

# Populate the newly inserted row 238 with the new weekly record.
$ws.Range("A238").Value = 7
$ws.Range("B238").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C238").Value = "Ñuble"
$ws.Range("D238").Value = 44461
$ws.Range("E238").Value = 16
$ws.Range("F238").Value = 100114014
$ws.Range("G238").Value = "Betarraga"
$ws.Range("H238").Value = "Sin especificar"
$ws.Range("I238").Value = "Primera"
$ws.Range("J238").Value = 300
$ws.Range("K238").Value = 700
$ws.Range("L238").Value = 750
$ws.Range("M238").Value = 725
$ws.Range("N238").Value = "$/paquete 5 unidades"
$ws.Range("O238").Value = "Región del Maule"
$ws.Range("P238").Value = 145
$ws.Range("Q238").Value = 5
$ws.Range("R238").Value = "Hortaliza"
